$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Purple color used throughout (RGB 0x70,0x30,0xA0 -> COM BGR-packed long)
$purple = 10498160

# --- Recolor existing "red" marker cells (fillId 3 / FFFF0000) to the
#     purple fill (fillId 4 / FF7030A0) -> matches style index 2 -> 3
$ws.Range("I10").Interior.Color = $purple
$ws.Range("J11").Interior.Color = $purple
$ws.Range("K12").Interior.Color = $purple
$ws.Range("K14").Interior.Color = $purple
$ws.Range("L18").Interior.Color = $purple
$ws.Range("L20").Interior.Color = $purple
$ws.Range("L21").Interior.Color = $purple
$ws.Range("L22").Interior.Color = $purple

# --- These two also flip to the new combined style: purple fill AND
#     purple font color (fontId 2 + fillId 4) -> matches style index 2 -> 5
$ws.Range("I23").Interior.Color = $purple
$ws.Range("I23").Font.Color = $purple
$ws.Range("J24").Interior.Color = $purple
$ws.Range("J24").Font.Color = $purple

# --- Insert 4 fresh rows before the old row 26 ("key" row); row 25 is
#     already a blank spacer row in the sheet so only 4 inserts are
#     needed to push the old row 26 down to row 30, matching the diff.
$ws.Rows("26:29").Insert()

# --- Populate the 5 new "Audio Assets" rows (25-29)
$ws.Range("A25").Value = "Glass Shatering"
$ws.Range("L25").Interior.Color = $purple
$ws.Range("L25").Font.Color = $purple

$ws.Range("A26").Value = "Walking Over Glass"
$ws.Range("L26").Interior.Color = $purple
$ws.Range("L26").Font.Color = $purple

$ws.Range("A27").Value = "Dumpster Sound"
$ws.Range("L27").Interior.Color = $purple

$ws.Range("A28").Value = "Falling Down"
$ws.Range("L28").Interior.Color = $purple

$ws.Range("A29").Value = "Broken UFO"
$ws.Range("L29").Interior.Color = $purple

# --- New cell added on the (now shifted) "key" row
$ws.Range("L30").Interior.Color = $purple

# --- Last row ("Scheduled") gains the combined purple fill+font style
$ws.Range("B32").Interior.Color = $purple
$ws.Range("B32").Font.Color = $purple

# --- Selection / view bookkeeping to mirror the saved workbook state
$ws.Range("Q18").Select()
